# Regenerate merged AHB files
# - Rename header columns from *_old / *_new to *_FV2404 / *_FV2410
# - Add an Excel Table (ListObject) over the data range
# - Freeze the header row (pane split)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row columns ---
$oldHeaders = @(
    "Segmentname_old","Segmentgruppe_old","Segment_old","Datenelement_old","Segment ID_old",
    "Code_old","Qualifier_old","Beschreibung_old","Bedingungsausdruck_old","Bedingung_old"
)
$fv2404Headers = @(
    "Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404","Segment ID_FV2404",
    "Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404","Bedingungsausdruck_FV2404","Bedingung_FV2404"
)
$newHeaders = @(
    "Segmentname_new","Segmentgruppe_new","Segment_new","Datenelement_new","Segment ID_new",
    "Code_new","Qualifier_new","Beschreibung_new","Bedingungsausdruck_new","Bedingung_new"
)
$fv2410Headers = @(
    "Segmentname_FV2410","Segmentgruppe_FV2410","Segment_FV2410","Datenelement_FV2410","Segment ID_FV2410",
    "Code_FV2410","Qualifier_FV2410","Beschreibung_FV2410","Bedingungsausdruck_FV2410","Bedingung_FV2410"
)

# Columns A-J hold the "_old" -> "_FV2404" headers
for ($i = 0; $i -lt $fv2404Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2404Headers[$i]
}

# Column K holds "diff" (unchanged)

# Columns L-U hold the "_new" -> "_FV2410" headers
for ($i = 0; $i -lt $fv2410Headers.Length; $i++) {
    $ws.Cells.Item(1, 12 + $i).Value = $fv2410Headers[$i]
}

# --- Add table (ListObject) over the data range ---
$tableRange = $ws.Range("A1:U76")
$list = $ws.ListObjects.Add(1, $tableRange, [System.Reflection.Missing]::Value, 1)
$list.Name = "Table1"

# --- Freeze header row ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
